$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price / Volume(1h) updates for existing rows ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.249.07"
$ws.Range("E2").Value = "  -1.59%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.750.59"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.61"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.91"
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.750.35"
$ws.Range("E7").Value = "  -2.29%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  -0.06%  "
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("E13").Value = "  +3.90%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.34"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.379.30"
$ws.Range("E15").Value = "  -2.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.753.70"
$ws.Range("E16").Value = "  -2.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.55"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.232.35"
$ws.Range("E18").Value = "  -1.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.17"
$ws.Range("E19").Value = "  -3.61%  "
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("E21").Value = "  -5.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.51"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("E23").Value = "  -2.91%  "
$ws.Range("E24").Value = "  -8.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.71"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("E27").Value = "  -0.79%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.31"
$ws.Range("E28").Value = "  +2.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.90"
$ws.Range("E30").Value = "  -2.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.900.74"
$ws.Range("E31").Value = "  -2.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.63"
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.40"
$ws.Range("E33").Value = "  -4.09%  "
$ws.Range("E34").Value = "  -4.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.09"
$ws.Range("E35").Value = "  -3.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.713.45"
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.80"
$ws.Range("E37").Value = "  +3.06%  "
$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("E39").Value = "  -2.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.994"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.310"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.68"
$ws.Range("E45").Value = "  -1.27%  "
$ws.Range("E46").Value = "  -2.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "397.03"
$ws.Range("E48").Value = "  -5.22%  "
$ws.Range("E49").Value = "  -8.42%  "

# --- Row 50/51: Monero inserted before VeChain, Arweave dropped off the list ---
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "138.79"
$ws.Range("E50").Value = "  -2.00%  "

$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0352"
$ws.Range("E51").Value = "  -2.50%  "
